$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header formatting (bold, border, alignment) from H1 into I1:J1,
# then set the new header labels.
$ws.Range("H1").Copy($ws.Range("I1:J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# I0 / IF data values for rows 2-82
$data = @(
    @(2,7,8),
    @(3,5,6),
    @(4,9,9),
    @(5,1,1),
    @(6,1,1),
    @(7,8,8),
    @(8,8,9),
    @(9,6,7),
    @(10,1,2),
    @(11,1,1),
    @(12,8,9),
    @(13,6,7),
    @(14,9,9),
    @(15,5,7),
    @(16,7,7),
    @(17,9,10),
    @(18,9,9),
    @(19,1,2),
    @(20,9,9),
    @(21,10,10),
    @(22,9,9),
    @(23,7,8),
    @(24,1,1),
    @(25,7,8),
    @(26,6,8),
    @(27,10,10),
    @(28,6,7),
    @(29,7,9),
    @(30,9,9),
    @(31,1,3),
    @(32,1,3),
    @(33,1,2),
    @(34,1,2),
    @(35,7,8),
    @(36,9,9),
    @(37,2,3),
    @(38,9,9),
    @(39,7,8),
    @(40,7,8),
    @(41,5,6),
    @(42,8,9),
    @(43,9,9),
    @(44,7,7),
    @(45,1,2),
    @(46,5,6),
    @(47,8,8),
    @(48,6,7),
    @(49,8,8),
    @(50,6,8),
    @(51,8,8),
    @(52,9,9),
    @(53,5,7),
    @(54,8,9),
    @(55,7,8),
    @(56,9,9),
    @(57,2,5),
    @(58,6,7),
    @(59,7,7),
    @(60,8,9),
    @(61,3,5),
    @(62,6,7),
    @(63,5,6),
    @(64,9,9),
    @(65,7,8),
    @(66,6,7),
    @(67,7,7),
    @(68,7,8),
    @(69,6,6),
    @(70,5,6),
    @(71,6,8),
    @(72,8,9),
    @(73,6,6),
    @(74,7,7),
    @(75,8,9),
    @(76,9,9),
    @(77,7,7),
    @(78,7,8),
    @(79,7,7),
    @(80,4,4),
    @(81,3,3),
    @(82,3,3)
)

foreach ($row in $data) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
